$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.807.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +7.32%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.807.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.43%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'250.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.40%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9991"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4951"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2809"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +8.08%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.48%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.806.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.35%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'17.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.92%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07095"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.6461"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.88%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.700"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.40%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'81.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.71%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'28.780.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +8.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.9990"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.06%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000007365"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.75%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.9987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +7.54%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.038.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.610"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.98%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.902"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.01%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.314"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.97%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'142.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.01%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'15.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +4.15%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.886"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +5.98%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'111.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.57%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.390"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.43%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.182"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.54%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08364"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.848"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.61%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04962"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +9.53%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.091"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.08%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.6704"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +7.63%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'HuobiToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'2.671"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.60%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.9518"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.98%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'RenderToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'2.155"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +5.44%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'MXToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'2.636"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +8.40%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01601"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.51%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.971"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.78%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.9985"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'100.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.33%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.4107"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.48%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'7.224"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.39%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.1226"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.72%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05489"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.87%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.116"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.00%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'31.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.44%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.308"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.3620"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.93%  "
$ws.Range("E51").Style = "Normal"
